# Scen_GRD_Constraints.xlsx - "Update af faste scenarie filer"
# Adds two new region columns (DE4, DE5) to the "DH data potentials" sheet,
# between the existing DE3 column and the "Name"/unit column, mirroring the
# existing DE1-DE3 layout (header, unit row, data rows, and the summary
# row at the bottom).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DH data potentials")

# --- Insert two new columns at I (pushes old I:M -> K:O) -------------------
$ws.Columns.Item(9).Insert()
$ws.Columns.Item(9).Insert()

# --- Header row (row 5): new region labels ---------------------------------
$ws.Range("I5").Value = "DE4"
$ws.Range("J5").Value = "DE5"

# --- Unit row (row 6): "MW" like the other region columns ------------------
$ws.Range("I6").Value = "MW"
$ws.Range("J6").Value = "MW"

# --- Row 7 (2010, existing) --------------------------------------------------
$ws.Range("F7").Formula = "=SUM(G7:J7)"
$ws.Range("G7").Value = 10
$ws.Range("H7").Formula = "=10"
$ws.Range("I7").Value = 10
$ws.Range("J7").Formula = "=G7"

# --- Row 8 (2012, existing) --------------------------------------------------
$ws.Range("F8").Formula = "=SUM(G8:J8)"
$ws.Range("G8").Value = 10
$ws.Range("H8").Formula = "=G8"
$ws.Range("I8").Formula = "=G8"
$ws.Range("J8").Formula = "=G8"

# --- Row 9 (2020, existing) --------------------------------------------------
$ws.Range("F9").Formula = "=SUM(G9:J9)"
$ws.Range("G9").Value = 10
$ws.Range("H9").Formula = "=10"
$ws.Range("I9").Formula = "=G9"
$ws.Range("J9").Formula = "=G9"

# --- Row 10 (2030, step 1) ---------------------------------------------------
$ws.Range("F10").Formula = "=SUM(G10:J10)"
$ws.Range("G10").Value = 222.3
$ws.Range("H10").Formula = "=G10"
$ws.Range("I10").Formula = "=G10"
$ws.Range("J10").Formula = "=G10"

# --- Row 11 (2050, step 1) ---------------------------------------------------
$ws.Range("F11").Formula = "=SUM(G11:J11)"
$ws.Range("G11").Formula = "=G10*2"
$ws.Range("H11").Formula = "=H10*2"
$ws.Range("I11").Formula = "=I10*2"
$ws.Range("J11").Formula = "=J10*2"

# --- Row 12 (2010, step 2) ---------------------------------------------------
$ws.Range("F12").Formula = "=SUM(G12:J12)"
$ws.Range("G12").Formula = "=G7/2"
$ws.Range("H12").Formula = "=H7/2"
$ws.Range("I12").Formula = "=H12"
$ws.Range("J12").Formula = "=G12"

# --- Row 13 (2012, step 2) ---------------------------------------------------
$ws.Range("F13").Formula = "=SUM(G13:J13)"
$ws.Range("G13").Formula = "=G12"
$ws.Range("H13").Formula = "=H8/2"
$ws.Range("I13").Formula = "=H13"
$ws.Range("J13").Formula = "=G13"

# --- Row 14 (2020, step 2) ---------------------------------------------------
$ws.Range("F14").Formula = "=SUM(G14:J14)"
$ws.Range("G14").Formula = "=G12"
$ws.Range("H14").Formula = "=H9/2"
$ws.Range("I14").Formula = "=H14"
$ws.Range("J14").Formula = "=G14"

# --- Row 15 (2030, step 2) ---------------------------------------------------
$ws.Range("F15").Formula = "=SUM(G15:J15)"
$ws.Range("G15").Formula = "=G10/2"
$ws.Range("H15").Formula = "=H10/2"
$ws.Range("I15").Formula = "=H15"
$ws.Range("J15").Formula = "=G15"

# --- Row 16 (2040, step 2) ---------------------------------------------------
$ws.Range("F16").Formula = "=SUM(G16:J16)"
$ws.Range("G16").Formula = "=G15"
$ws.Range("H16").Formula = "=H11/2"
$ws.Range("I16").Formula = "=H16"
$ws.Range("J16").Formula = "=G16"

# --- Row 17 (2010, transmission) ---------------------------------------------
$ws.Range("F17").Formula = "=SUM(G17:J17)"
$ws.Range("G17").Formula = "=G12/2"
$ws.Range("H17").Formula = "=H12/2"
$ws.Range("I17").Formula = "=I12/2"
$ws.Range("J17").Formula = "=J12/2"

# --- Row 18 (2012, transmission) ---------------------------------------------
$ws.Range("F18").Formula = "=SUM(G18:J18)"
$ws.Range("G18").Formula = "=G13/2"
$ws.Range("H18").Formula = "=H13/2"
$ws.Range("I18").Formula = "=H18"
$ws.Range("J18").Formula = "=G18"

# --- Row 19 (2020, transmission) ----------------------------------------------
$ws.Range("F19").Formula = "=SUM(G19:J19)"
$ws.Range("G19").Formula = "=G14/2"
$ws.Range("H19").Formula = "=H14/2"
$ws.Range("I19").Formula = "=H19"
$ws.Range("J19").Formula = "=G19"

# --- Row 20 (2030, transmission) ----------------------------------------------
$ws.Range("F20").Formula = "=SUM(G20:J20)"
$ws.Range("G20").Formula = "=G15/2"
$ws.Range("H20").Formula = "=H15/2"
$ws.Range("I20").Formula = "=H20"
$ws.Range("J20").Formula = "=G20"

# --- Row 21 (2040, transmission) ----------------------------------------------
$ws.Range("F21").Formula = "=SUM(G21:J21)"
$ws.Range("G21").Formula = "=G16/2"
$ws.Range("H21").Formula = "=H16/2"
$ws.Range("I21").Formula = "=H21"
$ws.Range("J21").Formula = "=G21"

# --- Row 22 (flat summary row) -----------------------------------------------
$ws.Range("I22").Value = 5
$ws.Range("J22").Value = 5
